$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 555, shifting existing rows 555:656 down to 556:657.
$ws.Rows(555).Insert()

# Populate the newly inserted row 555 with the new record's data.
$ws.Cells.Item(555, 1).Value = 9
$ws.Cells.Item(555, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(555, 3).Value = "Metropolitana"
$ws.Cells.Item(555, 4).Value = 45258
$ws.Cells.Item(555, 5).Value = 13
$ws.Cells.Item(555, 6).Value = 100112052
$ws.Cells.Item(555, 7).Value = "Albahaca"
$ws.Cells.Item(555, 8).Value = "Sin especificar"
$ws.Cells.Item(555, 9).Value = "Primera"
$ws.Cells.Item(555, 10).Value = 220
$ws.Cells.Item(555, 11).Value = 5000
$ws.Cells.Item(555, 12).Value = 6000
$ws.Cells.Item(555, 13).Value = 5545
$ws.Cells.Item(555, 14).Value = "$/docena de matas"
$ws.Cells.Item(555, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(555, 16).Value = 924
$ws.Cells.Item(555, 17).Value = 6
$ws.Cells.Item(555, 18).Value = "Hortaliza"

# Keep the date-formatted style (matching column D's other cells) on the new row's date cell.
$ws.Cells.Item(555, 4).NumberFormat = $ws.Cells.Item(556, 4).NumberFormat
